$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.02020935029036
$ws.Range("D2").Value = 1.028186279667556
$ws.Range("E2").Value = 1.021220136258265
$ws.Range("F2").Value = 1.035130770811574
$ws.Range("J2").Value = 1.025407753437642
$ws.Range("K2").Value = 1.03100382927611
$ws.Range("L2").Value = 1.024058110239223
$ws.Range("M2").Value = 1.037928256868433
$ws.Range("N2").Value = 1.012452158818634
$ws.Range("C3").Value = 1.02132920632262
$ws.Range("D3").Value = 1.029220324904257
$ws.Range("E3").Value = 1.02217679549443
$ws.Range("F3").Value = 1.036329240208038
$ws.Range("J3").Value = 1.026163808921131
$ws.Range("K3").Value = 1.031845146685901
$ws.Range("L3").Value = 1.02482076425516
$ws.Range("M3").Value = 1.038935022641683
$ws.Range("N3").Value = 1.012711453180477
$ws.Range("C4").Value = 1.022054374107675
$ws.Range("D4").Value = 1.02989025789283
$ws.Range("E4").Value = 1.022796603130266
$ws.Range("F4").Value = 1.037105892393236
$ws.Range("J4").Value = 1.026653067329834
$ws.Range("K4").Value = 1.032389772332395
$ws.Range("L4").Value = 1.025314422923699
$ws.Range("M4").Value = 1.039587063615902
$ws.Range("N4").Value = 1.012879013867214
$ws.Range("C5").Value = 1.022359366053031
$ws.Range("D5").Value = 1.03017209807532
$ws.Range("E5").Value = 1.023057358075439
$ws.Range("F5").Value = 1.037432675626427
$ws.Range("J5").Value = 1.026858761222235
$ws.Range("K5").Value = 1.032618789809561
$ws.Range("L5").Value = 1.02552199781915
$ws.Range("M5").Value = 1.03986132463169
$ws.Range("N5").Value = 1.012949403544407
$ws.Range("C6").Value = 1.022410583268608
$ws.Range("D6").Value = 1.030219432033327
$ws.Range("E6").Value = 1.023101150987446
$ws.Range("F6").Value = 1.037487560321049
$ws.Range("J6").Value = 1.026893298700969
$ws.Range("K6").Value = 1.032657246189083
$ws.Range("L6").Value = 1.025556852948643
$ws.Range("M6").Value = 1.039907382655752
$ws.Range("N6").Value = 1.012961219178527
$ws.Range("C7").Value = 1.022058448907758
$ws.Range("D7").Value = 1.029894023067784
$ws.Range("E7").Value = 1.022800086612803
$ws.Range("F7").Value = 1.037110257791373
$ws.Range("J7").Value = 1.026655815783625
$ws.Range("K7").Value = 1.032392832251587
$ws.Range("L7").Value = 1.025317196389368
$ws.Range("M7").Value = 1.039590727742831
$ws.Range("N7").Value = 1.012879954625311
$ws.Range("C8").Value = 1.020587698224991
$ws.Range("D8").Value = 1.028535566922472
$ws.Range("E8").Value = 1.02154328106708
$ws.Range("F8").Value = 1.035535558714177
$ws.Range("J8").Value = 1.02566325770935
$ws.Range("K8").Value = 1.031288107421598
$ws.Range("L8").Value = 1.024315817406213
$ws.Range("M8").Value = 1.03826837478237
$ws.Range("N8").Value = 1.012539834050499
$ws.Range("C9").Value = 1.018000199174955
$ws.Range("D9").Value = 1.026148197175385
$ws.Range("E9").Value = 1.019334648703178
$ws.Range("F9").Value = 1.032769622552556
$ws.Range("J9").Value = 1.023914545290447
$ws.Range("K9").Value = 1.029343251318956
$ws.Range("L9").Value = 1.022552566921314
$ws.Range("M9").Value = 1.035942774372515
$ws.Range("N9").Value = 1.01193881860297
$ws.Range("C10").Value = 1.016277946400031
$ws.Range("D10").Value = 1.024560915205688
$ws.Range("E10").Value = 1.017866280673537
$ws.Range("F10").Value = 1.030931612079898
$ws.Range("J10").Value = 1.022748931942025
$ws.Range("K10").Value = 1.02804788840283
$ws.Range("L10").Value = 1.021377945132637
$ws.Range("M10").Value = 1.034395417686236
$ws.Range("N10").Value = 1.011537020073727
$ws.Range("C11").Value = 1.015532832793117
$ws.Range("D11").Value = 1.023874621531023
$ws.Range("E11").Value = 1.017231423147907
$ws.Range("F11").Value = 1.030137137101096
$ws.Range("J11").Value = 1.022244252338499
$ws.Range("K11").Value = 1.027487265712157
$ws.Range("L11").Value = 1.020869528558607
$ws.Range("M11").Value = 1.033726112216235
$ws.Range("N11").Value = 1.011362772239941
$ws.Range("C12").Value = 1.015256158558723
$ws.Range("D12").Value = 1.02361985292866
$ws.Range("E12").Value = 1.016995752314614
$ws.Range("F12").Value = 1.029842242742463
$ws.Range("J12").Value = 1.022056797388697
$ws.Range("K12").Value = 1.02727906693356
$ws.Range("L12").Value = 1.020680710271986
$ws.Range("M12").Value = 1.033477608559682
$ws.Range("N12").Value = 1.011298008889888
$ws.Range("C13").Value = 1.015315501877299
$ws.Range("D13").Value = 1.023674494811127
$ws.Range("E13").Value = 1.017046298014399
$ws.Range("F13").Value = 1.029905489135307
$ws.Range("J13").Value = 1.022097006866215
$ws.Range("K13").Value = 1.027323724407619
$ws.Range("L13").Value = 1.020721211068086
$ws.Range("M13").Value = 1.033530908639986
$ws.Range("N13").Value = 1.011311902654673
$ws.Range("C14").Value = 1.015509960899997
$ws.Range("D14").Value = 1.023853559191552
$ws.Range("E14").Value = 1.017211939595774
$ws.Range("F14").Value = 1.030112756781235
$ws.Range("J14").Value = 1.022228757142771
$ws.Range("K14").Value = 1.027470055092403
$ws.Range("L14").Value = 1.020853920163529
$ws.Range("M14").Value = 1.033705568662979
$ws.Range("N14").Value = 1.011357419695646
$ws.Range("C15").Value = 1.015629785966006
$ws.Range("D15").Value = 1.023963906722703
$ws.Range("E15").Value = 1.017314015857491
$ws.Range("F15").Value = 1.030240488896373
$ws.Range("J15").Value = 1.022309933560282
$ws.Range("K15").Value = 1.027560219740043
$ws.Range("L15").Value = 1.020935690623516
$ws.Range("M15").Value = 1.033813196524736
$ws.Range("N15").Value = 1.011385458955995
$ws.Range("C16").Value = 1.016327410391851
$ws.Range("D16").Value = 1.02460648350824
$ws.Range("E16").Value = 1.017908434205112
$ws.Range("F16").Value = 1.030984368116399
$ws.Range("J16").Value = 1.022782426655572
$ws.Range("K16").Value = 1.02808510091451
$ws.Range("L16").Value = 1.021411691318446
$ws.Range("M16").Value = 1.034439852205829
$ws.Range("N16").Value = 1.011548578722275
$ws.Range("C17").Value = 1.016765180585314
$ws.Range("D17").Value = 1.025009825113517
$ws.Range("E17").Value = 1.018281552923018
$ws.Range("F17").Value = 1.031451357092458
$ws.Range("J17").Value = 1.02307881924232
$ws.Range("K17").Value = 1.028414419212292
$ws.Range("L17").Value = 1.021710328067462
$ws.Range("M17").Value = 1.034833126539842
$ws.Range("N17").Value = 1.011650828170646
$ws.Range("C18").Value = 1.017020585694692
$ws.Range("D18").Value = 1.025245184907689
$ws.Range("E18").Value = 1.018499278954748
$ws.Range("F18").Value = 1.031723878568846
$ws.Range("J18").Value = 1.023251703773305
$ws.Range("K18").Value = 1.028606531847662
$ws.Range("L18").Value = 1.021884537369873
$ws.Range("M18").Value = 1.035062585361173
$ws.Range("N18").Value = 1.01171044283295
$ws.Range("C19").Value = 1.017107682641319
$ws.Range("D19").Value = 1.02532545305531
$ws.Range("E19").Value = 1.018573533544294
$ws.Range("F19").Value = 1.031816824218712
$ws.Range("J19").Value = 1.023310653575343
$ws.Range("K19").Value = 1.028672041854562
$ws.Range("L19").Value = 1.021943941558649
$ws.Range("M19").Value = 1.035140836533131
$ws.Range("N19").Value = 1.011730765528751
$ws.Range("C20").Value = 1.016718205667346
$ws.Range("D20").Value = 1.02496654029778
$ws.Range("E20").Value = 1.018241511279631
$ws.Range("F20").Value = 1.03140123966609
$ws.Range("J20").Value = 1.023047018736897
$ws.Range("K20").Value = 1.028379083701165
$ws.Range("L20").Value = 1.021678285150477
$ws.Range("M20").Value = 1.034790924838147
$ws.Range("N20").Value = 1.011639860431314
$ws.Range("C21").Value = 1.015452694973948
$ws.Range("D21").Value = 1.023800825016871
$ws.Range("E21").Value = 1.017163158317031
$ws.Range("F21").Value = 1.03005171586077
$ws.Range("J21").Value = 1.022189959819832
$ws.Range("K21").Value = 1.027426963211209
$ws.Range("L21").Value = 1.020814839808249
$ws.Range("M21").Value = 1.033654132711969
$ws.Range("N21").Value = 1.011344017163625
$ws.Range("C22").Value = 1.014657562323476
$ws.Range("D22").Value = 1.023068769509993
$ws.Range("E22").Value = 1.016485985348941
$ws.Range("F22").Value = 1.029204425698047
$ws.Range("J22").Value = 1.021651124695448
$ws.Range("K22").Value = 1.026828566936166
$ws.Range("L22").Value = 1.020272132494018
$ws.Range("M22").Value = 1.032940000555417
$ws.Range("N22").Value = 1.011157777683048
$ws.Range("C23").Value = 1.015079025937757
$ws.Range("D23").Value = 1.02345676305327
$ws.Range("E23").Value = 1.01684488885698
$ws.Range("F23").Value = 1.029653475863673
$ws.Range("J23").Value = 1.021936768466199
$ws.Range("K23").Value = 1.027145765510301
$ws.Range("L23").Value = 1.020559815364791
$ws.Range("M23").Value = 1.033318517473021
$ws.Range("N23").Value = 1.011256528631875
$ws.Range("C24").Value = 1.016739431421465
$ws.Range("D24").Value = 1.024986098541521
$ws.Range("E24").Value = 1.018259604089706
$ws.Range("F24").Value = 1.031423885157302
$ws.Range("J24").Value = 1.023061388004713
$ws.Range("K24").Value = 1.028395050214526
$ws.Range("L24").Value = 1.021692763904883
$ws.Range("M24").Value = 1.034809993758544
$ws.Range("N24").Value = 1.01164481636003
$ws.Range("C25").Value = 1.018668642350705
$ws.Range("D25").Value = 1.026764631418412
$ws.Range("E25").Value = 1.019904919118999
$ws.Range("F25").Value = 1.033483633225869
$ws.Range("J25").Value = 1.02436659444327
$ws.Range("K25").Value = 1.029845829691816
$ws.Range("L25").Value = 1.023008254574207
$ws.Range("M25").Value = 1.036543460030293
$ws.Range("N25").Value = 1.012094393749036
